$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename headers to match new column naming scheme
$ws.Range("A1").Value = "Row ID"
$ws.Range("C1").Value = "Task"
$ws.Range("E1").Value = "Start Date"
$ws.Range("F1").Value = "End Date"

# Update selected/active cell in sheet view
$ws.Range("F1").Select()
